$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 47, shifting existing rows 47-54 down to 48-55
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new price-report entry
$ws.Range("A47").Value = 5
$ws.Range("B47").Value = "Macroferia Regional de Talca"
$ws.Range("C47").Value = "Maule"
$ws.Range("D47").Value = 44995
$ws.Range("E47").Value = 7
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100104
$ws.Range("H47").Value = "Frutos de pepita"
$ws.Range("I47").Value = 100104003
$ws.Range("J47").Value = "Membrillo"
$ws.Range("K47").Value = "Champion"
$ws.Range("L47").Value = "Especial"
$ws.Range("M47").Value = 540
$ws.Range("N47").Value = 12000
$ws.Range("O47").Value = 13000
$ws.Range("P47").Value = 12481
$ws.Range("Q47").Value = "$/caja 18 kilos granel"
$ws.Range("R47").Value = "Región de O'Higgins"
$ws.Range("S47").Value = 693
$ws.Range("T47").Value = 18
